$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1973.5
$ws.Range("J32").Value = 2340.1667
$ws.Range("L32").Value = 2340.1667
$ws.Range("N32").Value = -2992.1667

$ws.Range("H42").Value = 22
$ws.Range("I42").Value = 10
$ws.Range("J42").Value = 28
$ws.Range("K42").Value = 30
$ws.Range("L42").Value = 84
$ws.Range("M42").Value = 200
$ws.Range("N42").Value = -544

$ws.Range("H62").Value = 9941.286
$ws.Range("I62").Value = 8979.799999999999
$ws.Range("J62").Value = 12345
$ws.Range("K62").Value = 8979.799999999999
$ws.Range("L62").Value = 12345
$ws.Range("M62").Value = -8355.799999999999
$ws.Range("N62").Value = -13593

$ws.Range("H64").Value = 9181
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10496

$ws.Range("H65").Value = 9941.286
$ws.Range("I65").Value = 8979.799999999999
$ws.Range("J65").Value = 12345
$ws.Range("K65").Value = 44899
$ws.Range("L65").Value = 61725
$ws.Range("M65").Value = -41779
$ws.Range("N65").Value = -67965

$ws.Range("H67").Value = 9181
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11716

$ws.Range("H86").Value = 300103460
$ws.Range("I86").Value = 333336420
$ws.Range("J86").Value = 250253970
$ws.Range("K86").Value = 333336420
$ws.Range("L86").Value = 250253970
$ws.Range("M86").Value = -333335297
$ws.Range("N86").Value = -250256216

$ws.Range("H89").Value = 300103460
$ws.Range("I89").Value = 333336420
$ws.Range("J89").Value = 250253970
$ws.Range("K89").Value = 1666682100
$ws.Range("L89").Value = 1251269850
$ws.Range("M89").Value = -1666676484
$ws.Range("N89").Value = -1251281082

$ws.Range("H132").Value = 1875.7693
$ws.Range("I132").Value = 1875.7693
$ws.Range("K132").Value = 5627.3079
$ws.Range("M132").Value = -3097.3079

$ws.Range("H138").Value = 2513.61
$ws.Range("J138").Value = 2827.262
$ws.Range("L138").Value = 8481.786
$ws.Range("N138").Value = -18761.786

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 9000
$ws.Range("I18").Value = 9000
$ws.Range("K18").Value = 9000
$ws.Range("M18").Value = -8678

$ws.Range("H32").Value = 21042466
$ws.Range("I32").Value = 22728116
$ws.Range("K32").Value = 22728116
$ws.Range("M32").Value = -22727829

$ws.Range("H63").Value = 3997.4
$ws.Range("I63").Value = 2525
$ws.Range("J63").Value = 5469.8
$ws.Range("K63").Value = 2525
$ws.Range("L63").Value = 5469.8
$ws.Range("M63").Value = -1839
$ws.Range("N63").Value = -6841.8

$ws.Range("H66").Value = 3997.4
$ws.Range("I66").Value = 2525
$ws.Range("J66").Value = 5469.8
$ws.Range("K66").Value = 12625
$ws.Range("L66").Value = 27349
$ws.Range("M66").Value = -9193
$ws.Range("N66").Value = -34213

$ws.Range("H74").Value = 2030.5312
$ws.Range("I74").Value = 2135.682
$ws.Range("J74").Value = 1799.2
$ws.Range("K74").Value = 2135.682
$ws.Range("L74").Value = 1799.2
$ws.Range("M74").Value = -1261.682
$ws.Range("N74").Value = -3547.2

$ws.Range("H77").Value = 2030.5312
$ws.Range("I77").Value = 2135.682
$ws.Range("J77").Value = 1799.2
$ws.Range("K77").Value = 10678.41
$ws.Range("L77").Value = 8996
$ws.Range("M77").Value = -6310.41
$ws.Range("N77").Value = -17732

$ws.Range("H122").Value = 4742.256
$ws.Range("I122").Value = 4101.5186
$ws.Range("J122").Value = 5823.5
$ws.Range("K122").Value = 12304.5558
$ws.Range("L122").Value = 17470.5
$ws.Range("M122").Value = -9854.555800000002
$ws.Range("N122").Value = -22370.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2648
$ws.Range("I105").Value = 2115.4285
$ws.Range("K105").Value = 2115.4285
$ws.Range("M105").Value = -368.4285

$ws.Range("H107").Value = 3330.875
$ws.Range("I107").Value = 3037.6365
$ws.Range("K107").Value = 3037.6365
$ws.Range("M107").Value = -1117.6365

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 47.625
$ws.Range("I7").Value = 31.785715
$ws.Range("K7").Value = 31.785715
$ws.Range("M7").Value = 81.214285

$ws.Range("H16").Value = 1310.1818
$ws.Range("I16").Value = 1380.7
$ws.Range("K16").Value = 1380.7
$ws.Range("M16").Value = -1093.7

$ws.Range("H31").Value = 4582.7554
$ws.Range("I31").Value = 2276.9167
$ws.Range("K31").Value = 2276.9167
$ws.Range("M31").Value = -1981.9167

$ws.Range("H34").Value = 4582.7554
$ws.Range("I34").Value = 2276.9167
$ws.Range("K34").Value = 2276.9167
$ws.Range("M34").Value = -2074.9167

$ws.Range("H36").Value = 49999.668
$ws.Range("I36").Value = 49999.668
$ws.Range("K36").Value = 49999.668
$ws.Range("M36").Value = -49611.668

$ws.Range("H40").Value = 49999.668
$ws.Range("I40").Value = 49999.668
$ws.Range("K40").Value = 49999.668
$ws.Range("M40").Value = -49839.668

$ws.Range("H58").Value = 3399.2727
$ws.Range("I58").Value = 3030.2424
$ws.Range("K58").Value = 3030.2424
$ws.Range("M58").Value = -2827.2424

$ws.Range("H62").Value = 3649.5
$ws.Range("J62").Value = 4968
$ws.Range("L62").Value = 4968
$ws.Range("N62").Value = -6216

$ws.Range("H65").Value = 3649.5
$ws.Range("J65").Value = 4968
$ws.Range("L65").Value = 24840
$ws.Range("N65").Value = -31080

$ws.Range("H105").Value = 1959.2
$ws.Range("I105").Value = 1837.4615
$ws.Range("K105").Value = 1837.4615
$ws.Range("M105").Value = -90.46149999999989

$ws.Range("H113").Value = 1310.1818
$ws.Range("I113").Value = 1380.7
$ws.Range("K113").Value = 1380.7
$ws.Range("M113").Value = 789.3

$ws.Range("H132").Value = 5862.2
$ws.Range("I132").Value = 5862.2
$ws.Range("K132").Value = 17586.6
$ws.Range("M132").Value = -15056.6

$ws.Range("H136").Value = 3399.2727
$ws.Range("I136").Value = 3030.2424
$ws.Range("K136").Value = 9090.727200000001
$ws.Range("M136").Value = -6540.727200000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 4480.1333
$ws.Range("I133").Value = 4220.2
$ws.Range("K133").Value = 12660.6
$ws.Range("M133").Value = -7600.599999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 775000
$ws.Range("J139").Value = 775000
$ws.Range("L139").Value = 775000
$ws.Range("N139").Value = -785280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1516001.6
$ws.Range("I2").Value = 7500000
$ws.Range("J2").Value = 20002
$ws.Range("K2").Value = 7500000
$ws.Range("L2").Value = 20002
$ws.Range("M2").Value = -7499888
$ws.Range("N2").Value = -20226

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 39062.75
$ws.Range("J45").Value = 47083.668
$ws.Range("L45").Value = 47083.668
$ws.Range("N45").Value = -48065.668

$ws.Range("H62").Value = 4618.1816
$ws.Range("J62").Value = 6410
$ws.Range("L62").Value = 6410
$ws.Range("N62").Value = -7658

$ws.Range("H65").Value = 4618.1816
$ws.Range("J65").Value = 6410
$ws.Range("L65").Value = 32050
$ws.Range("N65").Value = -38290

$ws.Range("H113").Value = 629.5
$ws.Range("I113").Value = 471.14285
$ws.Range("K113").Value = 1413.42855
$ws.Range("M113").Value = 756.5714499999999

$ws.Range("H130").Value = 53994.5
$ws.Range("J130").Value = 53994.5
$ws.Range("L130").Value = 53994.5
$ws.Range("N130").Value = -64034.5
